$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing entry ---
# D2 becomes a true numeric cell (no text/string formatting) holding the new
# document number, replacing the old shared-string value.
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 4).Value = 6500118082

# F2 / G2 switch from raw numbers to formatted text amounts.
$ws.Cells.Item(2, 6).Value = " 1.00 "
$ws.Cells.Item(2, 7).Value = " 4,350.00 "

# --- Row 3: brand-new data-entry row (Alt+A "add new entry") ---
$ws.Cells.Item(3, 1).Value = "MBK"
$ws.Cells.Item(3, 2).Value = "01/10/68"
$ws.Cells.Item(3, 3).Value = "026959000"

$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 4).Value = 6500118083

$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "001"

$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = " 1.00 "

$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = " 22,495.82 "

# --- Sheet cosmetics to match the refreshed template ---
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Range("C14").Select()

# The sheet no longer carries the legacy "ignored error" suppressions now
# that D2/D3 are real numbers and the text cells are freshly (re)entered.
$ws.Cells.Item(2, 2).ClearIgnoredErrors()
$ws.Cells.Item(2, 3).ClearIgnoredErrors()
$ws.Cells.Item(2, 5).ClearIgnoredErrors()
